$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: clone cell formatting (number format / alignment / font) from existing
# rows so the new rows 78-80 pick up the same style indices used elsewhere in the
# sheet, without disturbing any existing cell.

# Row 78 formatting: A/C/D/E/G like row 74, B like row 75 (bold-ish "s=6" style)
$ws.Range("A74").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("B75").Copy()
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("C74").Copy()
$ws.Range("C78").PasteSpecial(-4122)
$ws.Range("D74").Copy()
$ws.Range("D78").PasteSpecial(-4122)
$ws.Range("E74").Copy()
$ws.Range("E78").PasteSpecial(-4122)
$ws.Range("G74").Copy()
$ws.Range("G78").PasteSpecial(-4122)

# Row 79 formatting: A:D like row 77
$ws.Range("A77:D77").Copy()
$ws.Range("A79:D79").PasteSpecial(-4122)

# Row 80 formatting: A:D like row 77
$ws.Range("A77:D77").Copy()
$ws.Range("A80:D80").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Step 2: fill in the new log entries (2014/8/20 entries from the changelog)

# Row 78 - 流程 / 改进 / 添加申请权限
$ws.Range("A78").Value = 41871
$ws.Range("B78").Value = "流程"
$ws.Range("C78").Value = "改进"
$ws.Range("D78").Value = "添加申请权限"
$ws.Range("E78").Value = "think_flow_type"
$ws.Range("G78").Value = "request_duy,report_duty"

# Row 79 - 邮件 / 改进 / 增加SSL链接功能
$ws.Range("A79").Value = 41871
$ws.Range("B79").Value = "邮件"
$ws.Range("C79").Value = "改进"
$ws.Range("D79").Value = "增加SSL链接功能"

# Row 80 - 附件 / 附件 / 修复不显示添加按钮的BUG
$ws.Range("A80").Value = 41871
$ws.Range("B80").Value = "附件"
$ws.Range("C80").Value = "附件"
$ws.Range("D80").Value = "修复不显示添加按钮的BUG"

# Match the explicit row height used throughout the rest of the log table
$ws.Rows.Item(78).RowHeight = 21
$ws.Rows.Item(79).RowHeight = 21
$ws.Rows.Item(80).RowHeight = 21

# --- Step 3: move the active selection the way the author left it
$ws.Range("E80").Select()
